# Trade #37 closed at 2026-02-17 13:23:34 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B5").Value = -0.9
$summary.Range("B6").Value = 37
$summary.Range("B9").Value = 40.54

# --- Strategy Status sheet ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 37
$status.Range("G4").Value = 40.54

# --- New trade row values shared by "All Trades" and "MarketMaking" sheets ---
$tradeNum = 37
$date = "2026-02-17"
$time = "13:23:28"
$strategy = "MarketMaking"
$side = "DOWN"
$entryPrice = 0.97
$exitPrice = 0.97
$status2 = "CLOSED"
$pnlPct = 0
$pnlDollar = 0
$capitalAfter = 98.34
$entrySlip = 0
$exitSlip = 0
$confidence = 0.6
$entryReason = "Normal spread capture: 19600 bps"
$exitReason = "early_exit"
$duration = 0.14

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 38

    $ws.Cells.Item($row, 1).Value = $tradeNum
    # Leading apostrophe forces this date-shaped string to stay literal text
    # (matches the inline string stored by the author) instead of being
    # auto-converted to a date serial number.
    $ws.Cells.Item($row, 2).Value = "'" + $date
    $ws.Cells.Item($row, 3).Value = $time
    $ws.Cells.Item($row, 4).Value = $strategy
    $ws.Cells.Item($row, 5).Value = $side
    $ws.Cells.Item($row, 6).Value = $entryPrice
    $ws.Cells.Item($row, 7).Value = $exitPrice
    $ws.Cells.Item($row, 8).Value = $status2
    $ws.Cells.Item($row, 9).Value = $pnlPct
    $ws.Cells.Item($row, 10).Value = $pnlDollar
    $ws.Cells.Item($row, 11).Value = $capitalAfter
    $ws.Cells.Item($row, 12).Value = $entrySlip
    $ws.Cells.Item($row, 13).Value = $exitSlip
    $ws.Cells.Item($row, 14).Value = $confidence
    $ws.Cells.Item($row, 15).Value = $entryReason
    $ws.Cells.Item($row, 16).Value = $exitReason
    $ws.Cells.Item($row, 17).Value = $duration
}
